# Applies the workbook edit described by the Lamia_Profits diff.
# Each Leve row (Table_<Job>) gets its recomputed market-board/profit
# columns (H..N) refreshed to the new scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 389
$ws.Range("I8").Value = 389
$ws.Range("K8").Value = 1167
$ws.Range("M8").Value = -1028

$ws.Range("H19").Value = 905.9
$ws.Range("I19").Value = 794.0769
$ws.Range("J19").Value = 1113.5714
$ws.Range("K19").Value = 794.0769
$ws.Range("L19").Value = 1113.5714
$ws.Range("M19").Value = -619.0769
$ws.Range("N19").Value = -1463.5714

$ws.Range("H33").Value = 477.6
$ws.Range("I33").Value = 415.35715
$ws.Range("J33").Value = 622.8333
$ws.Range("K33").Value = 415.35715
$ws.Range("L33").Value = 622.8333
$ws.Range("M33").Value = -186.35715
$ws.Range("N33").Value = -1080.8333

$ws.Range("H137").Value = 10640919
$ws.Range("I137").Value = 33335172
$ws.Range("J137").Value = 2988.5781
$ws.Range("K137").Value = 100005516
$ws.Range("L137").Value = 8965.7343
$ws.Range("M137").Value = -100002966
$ws.Range("N137").Value = -14065.7343

$ws.Range("H141").Value = 4222.933
$ws.Range("I141").Value = 4299.2856
$ws.Range("J141").Value = 4156.125
$ws.Range("K141").Value = 12897.8568
$ws.Range("L141").Value = 12468.375
$ws.Range("M141").Value = -7717.856800000001
$ws.Range("N141").Value = -22828.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7177.2905
$ws.Range("I32").Value = 5750.033
$ws.Range("K32").Value = 5750.033
$ws.Range("M32").Value = -5463.033

$ws.Range("H74").Value = 2279.5
$ws.Range("I74").Value = 2019.4117
$ws.Range("K74").Value = 2019.4117
$ws.Range("M74").Value = -1145.4117

$ws.Range("H77").Value = 2279.5
$ws.Range("I77").Value = 2019.4117
$ws.Range("K77").Value = 10097.0585
$ws.Range("M77").Value = -5729.058500000001

$ws.Range("H122").Value = 3851.1482
$ws.Range("I122").Value = 3481.6924
$ws.Range("K122").Value = 10445.0772
$ws.Range("M122").Value = -7995.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4428.3438
$ws.Range("I20").Value = 3820.625
$ws.Range("K20").Value = 3820.625
$ws.Range("M20").Value = -3573.625

$ws.Range("H86").Value = 3356.5217
$ws.Range("I86").Value = 2525.9473
$ws.Range("K86").Value = 2525.9473
$ws.Range("M86").Value = -1402.9473

$ws.Range("H89").Value = 3356.5217
$ws.Range("I89").Value = 2525.9473
$ws.Range("K89").Value = 12629.7365
$ws.Range("M89").Value = -7013.736499999999

$ws.Range("H99").Value = 1690
$ws.Range("I99").Value = 1593.6875
$ws.Range("J99").Value = 2075.25
$ws.Range("K99").Value = 1593.6875
$ws.Range("L99").Value = 2075.25
$ws.Range("M99").Value = -95.6875
$ws.Range("N99").Value = -5071.25

$ws.Range("H105").Value = 16705.375
$ws.Range("I105").Value = 16219.066
$ws.Range("K105").Value = 16219.066
$ws.Range("M105").Value = -14472.066

$ws.Range("H134").Value = 2400.359
$ws.Range("I134").Value = 2062.8857
$ws.Range("K134").Value = 6188.657099999999
$ws.Range("M134").Value = -3653.657099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45210.24
$ws.Range("I31").Value = 2261.3333
$ws.Range("J31").Value = 109633.6
$ws.Range("K31").Value = 2261.3333
$ws.Range("L31").Value = 109633.6
$ws.Range("M31").Value = -1966.3333
$ws.Range("N31").Value = -110223.6

$ws.Range("H34").Value = 45210.24
$ws.Range("I34").Value = 2261.3333
$ws.Range("J34").Value = 109633.6
$ws.Range("K34").Value = 2261.3333
$ws.Range("L34").Value = 109633.6
$ws.Range("M34").Value = -2059.3333
$ws.Range("N34").Value = -110037.6

$ws.Range("H99").Value = 2605.5806
$ws.Range("I99").Value = 2499.8333
$ws.Range("K99").Value = 2499.8333
$ws.Range("M99").Value = -1001.8333

$ws.Range("H126").Value = 2605.5806
$ws.Range("I126").Value = 2499.8333
$ws.Range("K126").Value = 7499.499899999999
$ws.Range("M126").Value = -5029.499899999999

$ws.Range("H132").Value = 2414.7727
$ws.Range("I132").Value = 2017.5238
$ws.Range("J132").Value = 10757
$ws.Range("K132").Value = 6052.5714
$ws.Range("L132").Value = 32271
$ws.Range("M132").Value = -3522.5714
$ws.Range("N132").Value = -37331

$ws.Range("H134").Value = 2169.7036
$ws.Range("I134").Value = 1675.6923
$ws.Range("K134").Value = 5027.0769
$ws.Range("M134").Value = -2492.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 5002
$ws.Range("J31").Value = 5002
$ws.Range("L31").Value = 15006
$ws.Range("N31").Value = -15582

$ws.Range("H34").Value = 4554.5454
$ws.Range("I34").Value = 5700
$ws.Range("J34").Value = 4300
$ws.Range("K34").Value = 17100
$ws.Range("L34").Value = 12900
$ws.Range("M34").Value = -17016
$ws.Range("N34").Value = -13068

$ws.Range("H56").Value = 5650.619
$ws.Range("I56").Value = 5650.619
$ws.Range("K56").Value = 5650.619
$ws.Range("M56").Value = -5120.619

$ws.Range("H107").Value = 453249.62
$ws.Range("I107").Value = 267.575
$ws.Range("J107").Value = 1078052.5
$ws.Range("K107").Value = 802.7249999999999
$ws.Range("L107").Value = 3234157.5
$ws.Range("M107").Value = 1117.275
$ws.Range("N107").Value = -3237997.5

$ws.Range("H113").Value = 2299
$ws.Range("J113").Value = 3578.4
$ws.Range("L113").Value = 10735.2
$ws.Range("N113").Value = -15075.2

$ws.Range("H132").Value = 4351
$ws.Range("I132").Value = 3435.8235
$ws.Range("J132").Value = 5906.8
$ws.Range("K132").Value = 30922.4115
$ws.Range("L132").Value = 53161.2
$ws.Range("M132").Value = -28392.4115
$ws.Range("N132").Value = -58221.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -11192

$ws.Range("H107").Value = 956.6667
$ws.Range("I107").Value = 983.3333
$ws.Range("K107").Value = 983.3333
$ws.Range("M107").Value = 936.6667

$ws.Range("H126").Value = 3454.842
$ws.Range("J126").Value = 3857.1538
$ws.Range("L126").Value = 11571.4614
$ws.Range("N126").Value = -16511.4614

$ws.Range("H132").Value = 3210.8865
$ws.Range("I132").Value = 2755.4324
$ws.Range("K132").Value = 8266.297200000001
$ws.Range("M132").Value = -5736.297200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5181.7856
$ws.Range("I22").Value = 1708.6
$ws.Range("K22").Value = 1708.6
$ws.Range("M22").Value = -1413.6

$ws.Range("H27").Value = 5181.7856
$ws.Range("I27").Value = 1708.6
$ws.Range("K27").Value = 1708.6
$ws.Range("M27").Value = -1601.6

$ws.Range("H61").Value = 3270.476
$ws.Range("I61").Value = 2483.75
$ws.Range("K61").Value = 2483.75
$ws.Range("M61").Value = -2281.75

$ws.Range("H113").Value = 3270.476
$ws.Range("I113").Value = 2483.75
$ws.Range("K113").Value = 2483.75
$ws.Range("M113").Value = -313.75

$ws.Range("H122").Value = 109931.16
$ws.Range("I122").Value = 133339.97
$ws.Range("J122").Value = 6263.5713
$ws.Range("K122").Value = 400019.91
$ws.Range("L122").Value = 18790.7139
$ws.Range("M122").Value = -397569.91
$ws.Range("N122").Value = -23690.7139

$ws.Range("H132").Value = 6261.037
$ws.Range("I132").Value = 4708.1763
$ws.Range("K132").Value = 14124.5289
$ws.Range("M132").Value = -11594.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6500.5
$ws.Range("I81").Value = 4000.2
$ws.Range("K81").Value = 8000.4
$ws.Range("M81").Value = -6939.4

$ws.Range("H84").Value = 6500.5
$ws.Range("I84").Value = 4000.2
$ws.Range("K84").Value = 40002
$ws.Range("M84").Value = -34698

$ws.Range("H113").Value = 420.64285
$ws.Range("I113").Value = 411.08694
$ws.Range("K113").Value = 1233.26082
$ws.Range("M113").Value = 936.73918

$ws.Range("H132").Value = 1877.2894
$ws.Range("I132").Value = 1370.0312
$ws.Range("J132").Value = 4582.6665
$ws.Range("K132").Value = 4110.0936
$ws.Range("L132").Value = 13747.9995
$ws.Range("M132").Value = -1580.0936
$ws.Range("N132").Value = -18807.9995

Write-Output "Applied 203 cell updates across 8 sheets"
